# "Bento object repository revisited"
#
# The FilesTab query (cell B4 on the "startup" sheet) drops the
# `File Type` and `Breed` columns from its RETURN clause, and the
# sheet's saved selection moves from D4 to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Whole Blood"] 
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $query

# Move the saved selection from D4 to C4.
$ws.Range("C4").Select()
